# Update 1.6.1: Added Remaining Material Analysis Report
#
# The "Материалы" (Materials) sheet used to ship with 33 pre-filled template
# rows (rows 2-34) that only held the F/H helper formulas (=D-E / =E-G,
# all evaluating to 0). Those placeholder rows are removed so the sheet is
# just the header row, ready for the new Remaining Material Analysis Report
# to drive its own rows. The view is also refreshed (zoom + scroll +
# selection) to match the new, much smaller sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Материалы" (sheet1.xml)
$ws.Activate()

# Remove the now-unused template rows 2:34 (shifts nothing up from below -
# there is no data past row 34), which also collapses the used range/
# dimension down to just the header row.
$ws.Range("A2:A34").EntireRow.Delete()

# Refresh the window: scroll so column D is the left-most visible column,
# zoom in to 190%, and leave the selection on J10.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 190
$ws.Range("J10").Select()
